$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update timestamp footer text
$ws.Range("A1").Value = "Datos actualizados a 16 de Junio de 2020 a las 15:52"

# Row 4
$ws.Range("B4").Value = 2183598
$ws.Range("C4").Value = 648
$ws.Range("E4").Value = 1175244
$ws.Range("G4").Value = 56
$ws.Range("H4").Value = 118339

# Row 7
$ws.Range("B7").Value = 344788
$ws.Range("C7").Value = 1762
$ws.Range("D7").Value = 181202
$ws.Range("E7").Value = 153660
$ws.Range("G7").Value = 11
$ws.Range("H7").Value = 9926

# Row 13
$ws.Range("B13").Value = 188086
$ws.Range("C13").Value = 42
$ws.Range("E13").Value = 6099
$ws.Range("G13").Value = 2
$ws.Range("H13").Value = 8887

# Row 19
$ws.Range("B19").Value = 136315
$ws.Range("C19").Value = 4267
$ws.Range("D19").Value = 89540
$ws.Range("E19").Value = 45723
$ws.Range("G19").Value = 41
$ws.Range("H19").Value = 1052

# Row 23
$ws.Range("B23").Value = 82077
$ws.Range("C23").Value = 1201
$ws.Range("D23").Value = 60461
$ws.Range("E23").Value = 21536
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = 80

# Row 27
$ws.Range("A27").Value = "Suecia"
$ws.Range("B27").Value = 53323
$ws.Range("C27").Value = 940
$ws.Range("D27").Value = 0
$ws.Range("E27").Value = 0
$ws.Range("G27").Value = 48
$ws.Range("H27").Value = 4939

# Row 28
$ws.Range("A28").Value = "Colombia"
$ws.Range("B28").Value = 53063
$ws.Range("D28").Value = 19952
$ws.Range("E28").Value = 31385
$ws.Range("H28").Value = 1726

# Row 32
$ws.Range("B32").Value = 42982
$ws.Range("C32").Value = 346
$ws.Range("D32").Value = 28861
$ws.Range("E32").Value = 13828
$ws.Range("G32").Value = 2
$ws.Range("H32").Value = 293

# Row 35
$ws.Range("B35").Value = 37336
$ws.Range("C35").Value = 300
$ws.Range("D35").Value = 23212
$ws.Range("E35").Value = 12602
$ws.Range("G35").Value = 2
$ws.Range("H35").Value = 1522

# Row 37
$ws.Range("D37").Value = 10174
$ws.Range("E37").Value = 21749
$ws.Range("G37").Value = 8
$ws.Range("H37").Value = 862

# Row 57
$ws.Range("B57").Value = 12426
$ws.Range("C57").Value = 59
$ws.Range("E57").Value = 659
$ws.Range("G57").Value = 1
$ws.Range("H57").Value = 256

# Row 58
$ws.Range("B58").Value = 12250
$ws.Range("C58").Value = 33
$ws.Range("D58").Value = 11125
$ws.Range("E58").Value = 527

# Row 69
$ws.Range("B69").Value = 8655
$ws.Range("C69").Value = 8
$ws.Range("E69").Value = 275

# Row 76
$ws.Range("B76").Value = 5328
$ws.Range("C76").Value = 65
$ws.Range("E76").Value = 1290

# Row 88
$ws.Range("B88").Value = 3860
$ws.Range("C88").Value = 133
$ws.Range("D88").Value = 1328
$ws.Range("E88").Value = 2427
$ws.Range("G88").Value = 1
$ws.Range("H88").Value = 105

# Row 93
$ws.Range("A93").Value = "Bosnia y Herzegovina"
$ws.Range("B93").Value = 3085
$ws.Range("C93").Value = 45
$ws.Range("D93").Value = 2178
$ws.Range("E93").Value = 739
$ws.Range("G93").Value = 3
$ws.Range("H93").Value = 168

# Row 94
$ws.Range("A94").Value = "Venezuela"
$ws.Range("B94").Value = 3062
$ws.Range("D94").Value = 835
$ws.Range("E94").Value = 2201
$ws.Range("H94").Value = 26

# Row 99
$ws.Range("B99").Value = 2273
$ws.Range("C99").Value = 11
$ws.Range("D99").Value = 1994
$ws.Range("E99").Value = 195

# Row 103
$ws.Range("B103").Value = 1914
$ws.Range("C103").Value = 9
$ws.Range("E103").Value = 532

# Row 106
$ws.Range("B106").Value = 1812
$ws.Range("C106").Value = 2
$ws.Range("E106").Value = 6

# Row 115
$ws.Range("B115").Value = 1473
$ws.Range("C115").Value = 9
$ws.Range("D115").Value = 889
$ws.Range("E115").Value = 552

# Row 158
$ws.Range("D158").Value = 325
$ws.Range("E158").Value = 9

# Row 206
$ws.Range("A206").Value = "Groenlandia"

# Row 207
$ws.Range("A207").Value = "Islas Malvinas"

# Row 208
$ws.Range("A208").Value = "Santa Sede"
$ws.Range("D208").Value = 12
$ws.Range("H208").Value = 0

# Row 209
$ws.Range("A209").Value = "Islas Turcas y Caicos"
$ws.Range("D209").Value = 11
$ws.Range("H209").Value = 1

# Row 213
$ws.Range("A213").Value = "Papua Nueva Guinea"
$ws.Range("D213").Value = 8
$ws.Range("H213").Value = 0

# Row 214
$ws.Range("A214").Value = "Islas Virgenes Britanicas"
$ws.Range("D214").Value = 7
$ws.Range("H214").Value = 1
